# Finished Week 13 logging
# Add new player "W.Smallwood" to the WR (Wide Receiver) sheet stats table,
# with all stat columns initialized to 0, and leave the WR sheet active
# with the selection resting on the next empty row (K10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WR")

# Make WR the active sheet (moves tabSelected from K to WR, and updates
# the workbook's activeTab bookView).
$ws.Activate()

# Log the new player's row: name in column A, stat columns B:J set to 0.
$ws.Range("A9").Value = "W.Smallwood"
$ws.Range("B9:J9").Value = 0

# Leave the selection on the next row down, ready for the next entry.
$ws.Range("K10").Select()
